# Lancers scrape run @ 2025-12-17 12:40:09 JST:
# - Prepend the newly discovered job posting as row 2 (pushing the nine
#   previously-seen rows down by one).
# - Refresh the "取得日時" (fetched-at) timestamp on every row to this run's
#   timestamp.
# - Column H ("スキル概要") grew one character wider this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-17 12:40:09"

# Shift the existing 8 job rows (old rows 2..9) down to rows 3..10, along
# with their values/styles, by inserting a fresh row above the old row 2.
$ws.Rows.Item(2).Insert()

# --- Fill in the newly inserted row with the new listing -----------------
$ws.Cells.Item(2, 1).Value = $newTimestamp
$ws.Cells.Item(2, 2).Value = "【急募】生成AI×業務効率化の実装を支援するエンジニア募集"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5444662"
$ws.Cells.Item(2, 7).Value = 385
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆効率化"

# --- Refresh the fetch timestamp on every row (now rows 2..10) -----------
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- Rebuild the F-column hyperlinks ---------------------------------
# Row insertion moves cell values/styles down automatically, but the
# worksheet's <hyperlinks> entries are anchored to fixed cell refs and do
# not shift with the rows, so every hyperlink (old + new) is rebuilt here,
# anchored on the URL text that now lives in each row.
$ws.Hyperlinks.Delete()
$urls = @(
    "https://www.lancers.jp/work/detail/5444662",
    "https://www.lancers.jp/work/detail/5450864",
    "https://www.lancers.jp/work/detail/5455415",
    "https://www.lancers.jp/work/detail/5455862",
    "https://www.lancers.jp/work/detail/5455513",
    "https://www.lancers.jp/work/detail/5455714",
    "https://www.lancers.jp/work/detail/5455675",
    "https://www.lancers.jp/work/detail/5016989",
    "https://www.lancers.jp/work/detail/5455422"
)
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $urls[$i]) | Out-Null
    # Hyperlinks.Add() stamps a fresh "applyFont" variant of the Hyperlink
    # cell style; bounce through Normal and back so every F-column cell
    # collapses back onto the single shared Hyperlink style (as before).
    $ws.Cells.Item($row, 6).Style = "Normal"
    $ws.Cells.Item($row, 6).Style = "Hyperlink"
}

# --- Column H widened from 12 to 13 characters ----------------------------
# ColumnWidth is expressed in points against the workbook's default font,
# which is offset from the raw "characters" unit stored in the xlsx by the
# Normal style's max digit width (5/6 of a character here); subtract that
# offset so the saved <col width> lands on exactly 13.
$ws.Columns.Item(8).ColumnWidth = 13 - (5 / 6)
